$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("H2").Value = 3.7
$ws.Range("J2").Value = 4.33
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 4.33
$ws.Range("Q2").Value = 1.73
$ws.Range("R2").Value = 2.1
$ws.Range("S2").Value = 2.75
$ws.Range("T2").Value = 1.44
$ws.Range("U2").Value = 1.33
$ws.Range("V2").Value = 3.25
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 13
$ws.Range("AC2").Value = 29
$ws.Range("AD2").Value = 34
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 7
$ws.Range("AI2").Value = 151
$ws.Range("AJ2").Value = 8.5
$ws.Range("AK2").Value = 9.5
$ws.Range("AN2").Value = 13

# Row 3
$ws.Range("J3").Value = 2.88
$ws.Range("O3").Value = 1.22
$ws.Range("P3").Value = 4.33
$ws.Range("Q3").Value = 1.8
$ws.Range("R3").Value = 2
$ws.Range("S3").Value = 2.75
$ws.Range("T3").Value = 1.44
$ws.Range("Y3").Value = 9.5

# Row 4
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 3.2
$ws.Range("I4").Value = 3
$ws.Range("J4").Value = 3.1
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 3.5
$ws.Range("N4").Value = 9.5
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 2.04
$ws.Range("R4").Value = 1.86
$ws.Range("S4").Value = 3.5
$ws.Range("T4").Value = 1.3
$ws.Range("W4").Value = 1.75
$ws.Range("X4").Value = 2
$ws.Range("Z4").Value = 12
$ws.Range("AA4").Value = 10
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 21
$ws.Range("AD4").Value = 29
$ws.Range("AE4").Value = 9.5
$ws.Range("AI4").Value = 201
$ws.Range("AJ4").Value = 9.5
$ws.Range("AK4").Value = 15
$ws.Range("AL4").Value = 11
$ws.Range("AM4").Value = 29
$ws.Range("AN4").Value = 23

# Row 6
$ws.Range("L6").Value = 3.2
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.53
$ws.Range("AA6").Value = 13
$ws.Range("AF6").Value = 6
$ws.Range("AR6").Value = 1.9
$ws.Range("AS6").Value = 2

# Row 10
$ws.Range("G10").Value = 2.35
$ws.Range("I10").Value = 3
$ws.Range("J10").Value = 3
$ws.Range("M10").Value = 1.06
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.3
$ws.Range("P10").Value = 3.4
$ws.Range("Q10").Value = 2.03
$ws.Range("R10").Value = 1.83
$ws.Range("S10").Value = 3.5
$ws.Range("T10").Value = 1.29
$ws.Range("AD10").Value = 26
$ws.Range("AJ10").Value = 10

# Row 17
$ws.Range("G17").Value = 2.15
$ws.Range("I17").Value = 3.8
$ws.Range("K17").Value = 1.91
$ws.Range("L17").Value = 4.5
$ws.Range("M17").Value = 1.11
$ws.Range("N17").Value = 6.5
$ws.Range("O17").Value = 1.5
$ws.Range("P17").Value = 2.63
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.5
$ws.Range("S17").Value = 5
$ws.Range("T17").Value = 1.17
$ws.Range("U17").Value = 1.57
$ws.Range("V17").Value = 2.25
$ws.Range("W17").Value = 2.1
$ws.Range("X17").Value = 1.67
$ws.Range("Z17").Value = 9
$ws.Range("AB17").Value = 19
$ws.Range("AD17").Value = 41
$ws.Range("AE17").Value = 6.5
$ws.Range("AG17").Value = 19
$ws.Range("AI17").Value = 501
$ws.Range("AR17").Value = 1.93
$ws.Range("AS17").Value = 1.93

# Row 30
$ws.Range("N30").Value = 8
$ws.Range("O30").Value = 1.4
$ws.Range("P30").Value = 2.75

# Row 32
$ws.Range("G32").Value = 3.15
$ws.Range("H32").Value = 3.6
$ws.Range("J32").Value = 3.55
$ws.Range("K32").Value = 2.22
$ws.Range("P32").Value = 3.45
$ws.Range("Q32").Value = 1.65
$ws.Range("R32").Value = 1.98
$ws.Range("Y32").Value = 11.5
$ws.Range("AC32").Value = 25
$ws.Range("AE32").Value = 12.5
$ws.Range("AF32").Value = 7.1
$ws.Range("AJ32").Value = 8.75
$ws.Range("AK32").Value = 10.5
$ws.Range("AM32").Value = 18.5

# Row 39
$ws.Range("G39").Value = 1.53
$ws.Range("H39").Value = 4.33
$ws.Range("I39").Value = 5.75
$ws.Range("J39").Value = 2.05
$ws.Range("K39").Value = 2.38
$ws.Range("M39").Value = 1.04
$ws.Range("N39").Value = 13
$ws.Range("O39").Value = 1.22
$ws.Range("P39").Value = 4
$ws.Range("Q39").Value = 1.7
$ws.Range("R39").Value = 2.1
$ws.Range("S39").Value = 2.75
$ws.Range("T39").Value = 1.4
$ws.Range("U39").Value = 1.33
$ws.Range("V39").Value = 3.25
$ws.Range("AE39").Value = 13
$ws.Range("AF39").Value = 8

# Row 40
$ws.Range("G40").Value = 3.75
$ws.Range("H40").Value = 3.2
$ws.Range("I40").Value = 2
$ws.Range("J40").Value = 4.33
$ws.Range("L40").Value = 2.75
$ws.Range("O40").Value = 1.33
$ws.Range("P40").Value = 3.25
$ws.Range("Q40").Value = 2.1
$ws.Range("R40").Value = 1.7
$ws.Range("U40").Value = 1.44
$ws.Range("V40").Value = 2.63
$ws.Range("W40").Value = 1.91
$ws.Range("X40").Value = 1.91
$ws.Range("Y40").Value = 10
$ws.Range("Z40").Value = 19
$ws.Range("AA40").Value = 13
$ws.Range("AC40").Value = 34
$ws.Range("AI40").Value = 301
$ws.Range("AK40").Value = 9
$ws.Range("AL40").Value = 9
$ws.Range("AM40").Value = 17
$ws.Range("AN40").Value = 17
$ws.Range("AO40").Value = 29

# Row 47
$ws.Range("G47").Value = 1.7
$ws.Range("H47").Value = 3.4
$ws.Range("I47").Value = 4.5
$ws.Range("L47").Value = 4.75
$ws.Range("M47").Value = 1.03
$ws.Range("N47").Value = 9.5
$ws.Range("Z47").Value = 8

# Row 49
$ws.Range("G49").Value = 2.8
$ws.Range("I49").Value = 2.63
$ws.Range("J49").Value = 3.6
$ws.Range("M49").Value = 1.08
$ws.Range("N49").Value = 8
$ws.Range("AK49").Value = 12

# Row 61
$ws.Range("G61").Value = 2.1
$ws.Range("I61").Value = 3.2
$ws.Range("L61").Value = 3.8
$ws.Range("AG61").Value = 13.5
$ws.Range("AH61").Value = 60
$ws.Range("AI61").Value = 450
$ws.Range("AJ61").Value = 10
$ws.Range("AK61").Value = 17

